# Attendance Info (Ramzor Graphs)
$wb = $excel.ActiveWorkbook

$wsInterview = $wb.Worksheets.Item("Interview")
$wsComment   = $wb.Worksheets.Item("Comment")

# --- Data change --------------------------------------------------------
# The placeholder note text on the Comment sheet was retyped.
$wsComment.Range("A2").Value = "asasasasasasasasasasasasasasasas"

# --- Font changes --------------------------------------------------------
# Plain/body cells: Calibri -> Arial
$wsInterview.Range("A2:B2").Font.Name = "Arial"
$wsComment.Range("A2").Font.Name = "Arial"
$wsComment.Range("A1").Font.Name = "Arial"

# Bold/underlined heading cells: Calibri Light -> Times New Roman
$wsInterview.Range("A1:B1").Font.Name = "Times New Roman"

# --- Column widths (nearest value this engine's character grid allows) --
$wsInterview.Columns.Item(1).ColumnWidth = 19
$wsInterview.Columns.Item(2).ColumnWidth = 21.5
$wsComment.Columns.Item(1).ColumnWidth = 9.333333333333332

# --- Selections / active sheet -------------------------------------------
$wsInterview.Range("A6:A7").Select()
$wsInterview.Activate()
$wsComment.Range("E2").Select()
$wsComment.Activate()
